# Automated map update: remove resolved/stale "Caso" entries from the INCO sheet.
# This mirrors the upstream automation's commit
# "Actualización automática del mapa" which deletes the rows whose "Caso"
# (column A) matches one of the IDs below, shifting all subsequent rows up.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("INCO")

$casosToRemove = @("6180", "6580", "6486", "6695", "6478")

# Find the last used row/column so we know how far to scan.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row   # xlUp = -4162

foreach ($caso in $casosToRemove) {
    for ($r = $lastRow; $r -ge 2; $r--) {
        $val = $ws.Cells.Item($r, 1).Value()
        if ($null -ne $val -and "$val" -eq $caso) {
            $ws.Rows.Item($r).Delete() | Out-Null
            break
        }
    }
    $lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
}
